$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell H1 (bold, bordered, centered
# style) onto the two new header cells before setting their text, so I1/J1
# end up sharing the same cell style as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I and J columns for each data row (2-16): I = 1, J = same value as H
for ($r = 2; $r -le 16; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
